$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1: bump the "Last status check on" timestamp from 01:15 to 01:30
$ws.Range("F1").Value = "Last status check on: 26.02.2022 01:30"

# D5: was stored as inline string "+0.4" -> becomes a real numeric value 0.4
$ws.Range("D5").Value = 0.4

# E5: was stored as inline string "2022-02-26 01:17:15" -> becomes a real
# numeric Excel date/time serial value, formatted like the other Old Datum cells
$ws.Range("E5").Value = 44618.05364583333
$ws.Range("E5").NumberFormat = "YYYY-MM-DD HH:MM:SS"
